$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text number format so numeric-looking price strings (e.g. "1.00")
# are preserved verbatim as text rather than being parsed into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '48.814.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.640.30'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '322.76'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.540'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.52'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.84'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0809'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.21'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.045.16'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.633.53'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.866'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '48.832.76'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.82'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.67'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.88'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.76%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.51'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.29'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.08'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.13'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.52%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.95'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.137'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.46'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.25'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0797'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.74%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.16'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.29'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.74'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.70%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.15'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.37%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.067.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.77%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.12'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.83%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.97'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '58.60'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.16'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.52%  '
